# daily auto push: 2026-01-14 18:49 UTC
# Two new schedule rows (2026/01/14 and 2026/01/15) are inserted into the
# rolling schedule block, shifting the existing rows 639:680 down to 641:682.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("639:640").Insert()

# Row 639: 2026/01/14 (Wed)
$ws.Range("A639").Value = "'2026/01/14"
$ws.Range("A639").ClearFormats()
$ws.Range("B639").Value = "水"
$ws.Range("C639").Value = 23
$ws.Range("D639").Value = 201

# Row 640: 2026/01/15 (Thu)
$ws.Range("A640").Value = "'2026/01/15"
$ws.Range("A640").ClearFormats()
$ws.Range("B640").Value = "木"
$ws.Range("C640").Value = 2
$ws.Range("D640").Value = 201
